# taskset creation and optimization algo beginning
#
# Updates the benchmark numbers recorded in sheets 1-4, renames sheet 4
# ("ro_CCL-code_CCM" -> "ro_CCM-code_CCM", fixing a typo), and adds a new
# sheet 5 ("ro_RAM-code_CCM") holding the RAM/CCM taskset results.

$wb = $excel.ActiveWorkbook

# --- sheet 1: ro_FLASH-code_FLASH ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = 15948
$ws1.Range("C2").Value = 22051
$ws1.Range("D2").Value = 27695
$ws1.Range("B3").Value = 0.20262
$ws1.Range("C3").Value = 0.16208
$ws1.Range("D3").Value = 0.18926
$ws1.Range("B4").Value = "(4937.4, 5140.02)"
$ws1.Range("C4").Value = "(8115.6, 8277.68)"
$ws1.Range("D4").Value = "(10808.46, 10997.72)"
$ws1.Range("B5").Value = 10.664
$ws1.Range("C5").Value = 11.794
$ws1.Range("D5").Value = 17.297

# --- sheet 2: ro_FLASH-code_CCM ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = 13689
$ws2.Range("C2").Value = 23286
$ws2.Range("D2").Value = 31922
$ws2.Range("B3").Value = 0.20258
$ws2.Range("C3").Value = 0.12178
$ws2.Range("D3").Value = 0.0948
$ws2.Range("B4").Value = "(5336.42, 5539.0)"
$ws2.Range("C4").Value = "(8474.0, 8595.78)"
$ws2.Range("D4").Value = "(11194.08, 11288.88)"
$ws2.Range("B5").Value = 9.151
$ws2.Range("C5").Value = 9.358000000000001
$ws2.Range("D5").Value = 9.986000000000001

# --- sheet 3: ro_CCM-code_FLASH ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = 13989
$ws3.Range("C2").Value = 23348
$ws3.Range("D2").Value = 31420
$ws3.Range("B3").Value = 0.2027
$ws3.Range("C3").Value = 0.12186
$ws3.Range("D3").Value = 0.09497999999999999
$ws3.Range("B4").Value = "(5735.34, 5938.04)"
$ws3.Range("C4").Value = "(8792.12, 8913.98)"
$ws3.Range("D4").Value = "(11485.3, 11580.28)"
$ws3.Range("B5").Value = 9.356999999999999
$ws3.Range("C5").Value = 9.388999999999999
$ws3.Range("D5").Value = 9.848000000000001

# --- sheet 4: renamed ro_CCL-code_CCM -> ro_CCM-code_CCM ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "ro_CCM-code_CCM"
$ws4.Range("B2").Value = 11837
$ws4.Range("C2").Value = 22800
$ws4.Range("D2").Value = 33660
$ws4.Range("B3").Value = 0.20258
$ws4.Range("C3").Value = 0.1015
$ws4.Range("D3").Value = 0.0678
$ws4.Range("B4").Value = "(6134.36, 6336.94)"
$ws4.Range("C4").Value = "(9110.26, 9211.76)"
$ws4.Range("D4").Value = "(11776.66, 11844.46)"
$ws4.Range("B5").Value = 7.913
$ws4.Range("C5").Value = 7.637
$ws4.Range("D5").Value = 7.531

# --- new sheet 5: ro_RAM-code_CCM (inserted right after sheet 4) ---
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "ro_RAM-code_CCM"

# Header row ("24"/"48"/"72") must stay TEXT, not be coerced to numbers.
$ws5.Range("B1:D1").NumberFormat = "@"
$ws5.Range("B1").Value = "24"
$ws5.Range("C1").Value = "48"
$ws5.Range("D1").Value = "72"

$ws5.Range("A2").Value = "intensity"
$ws5.Range("B2").Value = 11434
$ws5.Range("C2").Value = 22037
$ws5.Range("D2").Value = 32561

$ws5.Range("A3").Value = "runtime"
$ws5.Range("B3").Value = 0.20258
$ws5.Range("C3").Value = 0.10154
$ws5.Range("D3").Value = 0.06802

$ws5.Range("A4").Value = "timestamp"
$ws5.Range("B4").Value = "(6932.28, 7134.86)"
$ws5.Range("C4").Value = "(9725.94, 9827.48)"
$ws5.Range("D4").Value = "(12331.96, 12399.98)"

$ws5.Range("A5").Value = "energy"
$ws5.Range("B5").Value = 7.644
$ws5.Range("C5").Value = 7.384
$ws5.Range("D5").Value = 7.309

# Match the bold/centered/bordered formatting used by sheets 2-4's header
# row and label column by copying the existing formatting over.
$ws4.Range("B1:D1").Copy()
$ws5.Range("B1:D1").PasteSpecial(-4122)

$ws4.Range("A2:A5").Copy()
$ws5.Range("A2:A5").PasteSpecial(-4122)

# Match the page margins used by sheets 2-4 (0.75"/1"/0.5").
$ws5.PageSetup.LeftMargin = 54
$ws5.PageSetup.RightMargin = 54
$ws5.PageSetup.TopMargin = 72
$ws5.PageSetup.BottomMargin = 72
$ws5.PageSetup.HeaderMargin = 36
$ws5.PageSetup.FooterMargin = 36

# Restore sheet 1 as the active/selected tab (unchanged by this edit).
$ws1.Activate()
[void]$ws1.Range("A1").Select()
